$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 107; existing rows 107+ shift down to 108+
$ws.Rows.Item(107).Insert()

# Populate the newly inserted row 107 with the new record
$ws.Cells.Item(107, 1).Value = 1
$ws.Cells.Item(107, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(107, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(107, 4).Value = 44977
$ws.Cells.Item(107, 5).Value = 15
$ws.Cells.Item(107, 6).Value = "Fruta"
$ws.Cells.Item(107, 7).Value = 100106
$ws.Cells.Item(107, 8).Value = "Oleaginosos"
$ws.Cells.Item(107, 9).Value = 100106002
$ws.Cells.Item(107, 10).Value = "Palta"
$ws.Cells.Item(107, 11).Value = "Hass"
$ws.Cells.Item(107, 12).Value = "Segunda"
$ws.Cells.Item(107, 13).Value = 200
$ws.Cells.Item(107, 14).Value = 25000
$ws.Cells.Item(107, 15).Value = 26000
$ws.Cells.Item(107, 16).Value = 25500
$ws.Cells.Item(107, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(107, 18).Value = "Perú"
$ws.Cells.Item(107, 19).Value = 2550
$ws.Cells.Item(107, 20).Value = 10
